$wb = $excel.ActiveWorkbook

# This script applies updated currentAveragePrice / profit calculations
# (columns H-N) produced by the scheduled market-data refresh runner,
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets.

$ws = $wb.Worksheets.Item("ALC")
# Row 2: Mercury Rising | Quicksilver
$ws.Range("H2").Value = 84.36364
$ws.Range("I2").Value = 88.71429000000001
$ws.Range("J2").Value = 76.75
$ws.Range("K2").Value = 88.71429000000001
$ws.Range("L2").Value = 76.75
$ws.Range("M2").Value = 24.28570999999999
$ws.Range("N2").Value = -302.75
# Row 12: Don't Be So Tallow | Beeswax
$ws.Range("H12").Value = 255.125
$ws.Range("I12").Value = 117.8
$ws.Range("K12").Value = 117.8
$ws.Range("M12").Value = 52.2
# Row 43: Growing Is Knowing | Growth Formula Gamma
$ws.Range("H43").Value = 3857403.8
$ws.Range("J43").Value = 5000
$ws.Range("L43").Value = 5000
$ws.Range("N43").Value = -5138
# Row 58: A Matter of Vital Importance | Mega-Potion of Vitality
$ws.Range("H58").Value = 3390.5
$ws.Range("J58").Value = 8899.666999999999
$ws.Range("L58").Value = 26699.001
$ws.Range("N58").Value = -26999.001
# Row 132: Fast-forwarding Flora | Growth Formula Lambda
$ws.Range("H132").Value = 115663.27
$ws.Range("I132").Value = 289453.88
$ws.Range("K132").Value = 868361.64
$ws.Range("M132").Value = -865831.64
# Row 137: Cutting Edge of Culinary Quality | Magnesia Whetstone
$ws.Range("H137").Value = 2521.0557
$ws.Range("J137").Value = 4167.1665
$ws.Range("L137").Value = 12501.4995
$ws.Range("N137").Value = -17601.4995
# Row 138: All-night Crafting | Cunning Craftsman's Tisane
$ws.Range("H138").Value = 5501.2
$ws.Range("I138").Value = 1964.9333
$ws.Range("K138").Value = 5894.7999
$ws.Range("M138").Value = -754.7999

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust | Steel Ingot
$ws.Range("H32").Value = 2649.8235
$ws.Range("I32").Value = 2088.6304
$ws.Range("K32").Value = 2088.6304
$ws.Range("M32").Value = -1801.6304
# Row 122: Haste for High Durium | High Durium Nugget
$ws.Range("H122").Value = 4873.478
$ws.Range("J122").Value = 6902.8184
$ws.Range("L122").Value = 20708.4552
$ws.Range("N122").Value = -25608.4552
# Row 140: A Hand for a Deckhand | Ra'Kaznar Gloves of Scouting
$ws.Range("H140").Value = 114999.5
$ws.Range("J140").Value = 114999.5
$ws.Range("L140").Value = 114999.5
$ws.Range("N140").Value = -125359.5
# Row 141: Essays on Equipment | Ra'Kaznar Greaves of Maiming
$ws.Range("H141").Value = 134482
$ws.Range("J141").Value = 139476
$ws.Range("L141").Value = 139476
$ws.Range("N141").Value = -149836

$ws = $wb.Worksheets.Item("BSM")
# Row 105: Ingot to Wing It | Molybdenum Ingot
$ws.Range("H105").Value = 150001580
$ws.Range("I105").Value = 150001580
$ws.Range("K105").Value = 150001580
$ws.Range("M105").Value = -149999833
# Row 107: The Gold Experience | Deepgold Nugget
$ws.Range("H107").Value = 1484.591
$ws.Range("I107").Value = 1401.5714
$ws.Range("K107").Value = 1401.5714
$ws.Range("M107").Value = 518.4286
# Row 134: Ruthenium Supremium | Ruthenium Ingot
$ws.Range("H134").Value = 5894.0713
$ws.Range("I134").Value = 4489
$ws.Range("K134").Value = 13467
$ws.Range("M134").Value = -10932

$ws = $wb.Worksheets.Item("CRP")
# Row 7: Gridania's Got Talent | Maple Lumber
$ws.Range("H7").Value = 117.958336
$ws.Range("I7").Value = 67.59999999999999
$ws.Range("J7").Value = 201.88889
$ws.Range("K7").Value = 67.59999999999999
$ws.Range("L7").Value = 201.88889
$ws.Range("M7").Value = 45.40000000000001
$ws.Range("N7").Value = -427.88889
# Row 15: On the Move | Ragstone Grinding Wheel
$ws.Range("H15").Value = 2833.6667
$ws.Range("I15").Value = 1250.5
$ws.Range("K15").Value = 1250.5
$ws.Range("M15").Value = -1080.5
# Row 31: Wall Not Found | Walnut Lumber
$ws.Range("H31").Value = 33335560
$ws.Range("J31").Value = 6830.6665
$ws.Range("L31").Value = 6830.6665
$ws.Range("N31").Value = -7420.6665
# Row 34: Armoires of the Rich and Famous | Walnut Lumber
$ws.Range("H34").Value = 33335560
$ws.Range("J34").Value = 6830.6665
$ws.Range("L34").Value = 6830.6665
$ws.Range("N34").Value = -7234.6665
# Row 99: O Pine | Pine Lumber
$ws.Range("H99").Value = 13117.833
$ws.Range("I99").Value = 22084.5
$ws.Range("J99").Value = 8634.5
$ws.Range("K99").Value = 22084.5
$ws.Range("L99").Value = 8634.5
$ws.Range("M99").Value = -20586.5
$ws.Range("N99").Value = -11630.5
# Row 108: Just Starting Out | White Oak Fishing Rod
$ws.Range("H108").Value = 53000
$ws.Range("J108").Value = 53000
$ws.Range("L108").Value = 53000
$ws.Range("N108").Value = -60680
# Row 126: A Better Conductor | Red Pine Lumber
$ws.Range("H126").Value = 13117.833
$ws.Range("I126").Value = 22084.5
$ws.Range("J126").Value = 8634.5
$ws.Range("K126").Value = 66253.5
$ws.Range("L126").Value = 25903.5
$ws.Range("M126").Value = -63783.5
$ws.Range("N126").Value = -30843.5
# Row 134: Wood You Be Quiet | Ceiba Lumber
$ws.Range("H134").Value = 1376.3077
$ws.Range("I134").Value = 1376.3077
$ws.Range("K134").Value = 4128.9231
$ws.Range("M134").Value = -1593.9231

$ws = $wb.Worksheets.Item("CUL")
# Row 5: What a Sap | Maple Syrup
$ws.Range("H5").Value = 901.7059
$ws.Range("J5").Value = 1364.2858
$ws.Range("L5").Value = 4092.8574
$ws.Range("N5").Value = -4316.857400000001
# Row 7: It's Always Sunny in Vylbrand | Raisins
$ws.Range("H7").Value = 400
$ws.Range("I7").Value = 400
$ws.Range("K7").Value = 1200
$ws.Range("M7").Value = -1088
# Row 55: Pagan Pastries | Pastry Fish
$ws.Range("H55").Value = 4406.25
$ws.Range("J55").Value = 4884.615
$ws.Range("L55").Value = 14653.845
$ws.Range("N55").Value = -15007.845
# Row 119: Super Dark Times | Risotto al Nero
$ws.Range("H119").Value = 13571.286
$ws.Range("I119").Value = 5666.6665
$ws.Range("K119").Value = 16999.9995
$ws.Range("M119").Value = -12161.9995
# Row 122: Salt of the North | Northern Sea Salt
$ws.Range("H122").Value = 859.8
$ws.Range("I122").Value = 699.75
$ws.Range("K122").Value = 6297.75
$ws.Range("M122").Value = -3847.75
# Row 126: Imperial Palate | Glory Be Soup
$ws.Range("H126").Value = 15038.777
$ws.Range("I126").Value = 12192.714
$ws.Range("K126").Value = 36578.142
$ws.Range("M126").Value = -31638.142
# Row 128: A Historical Flavor | Skyr
$ws.Range("H128").Value = 339073
$ws.Range("I128").Value = 339073
$ws.Range("K128").Value = 1017219
$ws.Range("M128").Value = -1012239
# Row 131: The Mountain Steeped | Tsai tou Vounou
$ws.Range("H131").Value = 18187202
$ws.Range("J131").Value = 7719320
$ws.Range("L131").Value = 23157960
$ws.Range("N131").Value = -23168040
# Row 135: Not-so-secret Ingredient | Royal Maple Syrup
$ws.Range("H135").Value = 901.7059
$ws.Range("J135").Value = 1364.2858
$ws.Range("L135").Value = 12278.5722
$ws.Range("N135").Value = -17348.5722

$ws = $wb.Worksheets.Item("GSM")
# Row 102: Put the Metal to the Peddle | Durium Ingot
$ws.Range("H102").Value = 22182234
$ws.Range("I102").Value = 30007586
$ws.Range("K102").Value = 30007586
$ws.Range("M102").Value = -30005964
# Row 122: Awarding Academic Excellence | Ametrine
$ws.Range("H122").Value = 616788.0600000001
$ws.Range("I122").Value = 1225848.1
$ws.Range("J122").Value = 7728
$ws.Range("K122").Value = 3677544.3
$ws.Range("L122").Value = 23184
$ws.Range("M122").Value = -3675094.3
$ws.Range("N122").Value = -28084
# Row 126: Gold Rush Order | Phrygian Gold Ingot
$ws.Range("H126").Value = 3433.8125
$ws.Range("I126").Value = 2174.8
$ws.Range("J126").Value = 7930.2856
$ws.Range("K126").Value = 6524.400000000001
$ws.Range("L126").Value = 23790.8568
$ws.Range("M126").Value = -4054.400000000001
$ws.Range("N126").Value = -28730.8568
# Row 132: On Board for Lar | Lar Ingot
$ws.Range("H132").Value = 3470.3794
$ws.Range("I132").Value = 2839.9524
$ws.Range("K132").Value = 8519.8572
$ws.Range("M132").Value = -5989.8572

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban | Leather
$ws.Range("H7").Value = 4673.28
$ws.Range("I7").Value = 3812.4285
$ws.Range("J7").Value = 5768.909
$ws.Range("K7").Value = 3812.4285
$ws.Range("L7").Value = 5768.909
$ws.Range("M7").Value = -3700.4285
$ws.Range("N7").Value = -5992.909
# Row 61: Spelling Me Softly | Raptor Leather
$ws.Range("H61").Value = 3156.8333
$ws.Range("I61").Value = 3156.8333
$ws.Range("K61").Value = 3156.8333
$ws.Range("M61").Value = -2954.8333
# Row 113: Peace in Rest | Atrociraptor Leather
$ws.Range("H113").Value = 3156.8333
$ws.Range("I113").Value = 3156.8333
$ws.Range("K113").Value = 3156.8333
$ws.Range("M113").Value = -986.8332999999998
# Row 126: Battered Books | Saiga Leather
$ws.Range("H126").Value = 4673.28
$ws.Range("I126").Value = 3812.4285
$ws.Range("J126").Value = 5768.909
$ws.Range("K126").Value = 11437.2855
$ws.Range("L126").Value = 17306.727
$ws.Range("M126").Value = -8967.2855
$ws.Range("N126").Value = -22246.727
# Row 132: Tenets of Tanning | Silver Lobo Leather
$ws.Range("H132").Value = 3690.96
$ws.Range("I132").Value = 3069.3
$ws.Range("J132").Value = 4623.45
$ws.Range("K132").Value = 9207.900000000001
$ws.Range("L132").Value = 13870.35
$ws.Range("M132").Value = -6677.900000000001
$ws.Range("N132").Value = -18930.35

$ws = $wb.Worksheets.Item("WVR")
# Row 107: Flax Wax | Bright Linen Yarn
$ws.Range("H107").Value = 1515.7142
$ws.Range("I107").Value = 1641.3684
$ws.Range("K107").Value = 4924.1052
$ws.Range("M107").Value = -3004.1052
# Row 126: A Polished Purchase | Snow Linen
$ws.Range("H126").Value = 1993.8
$ws.Range("I126").Value = 1494.75
$ws.Range("J126").Value = 3990
$ws.Range("K126").Value = 4484.25
$ws.Range("L126").Value = 11970
$ws.Range("M126").Value = -2014.25
$ws.Range("N126").Value = -16910
# Row 132: Comfy Cabins | Snow Cotton Cloth
$ws.Range("H132").Value = 10105376
$ws.Range("I132").Value = 1161504
$ws.Range("K132").Value = 3484512
$ws.Range("M132").Value = -3481982

